$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as plain text, forcing Excel to keep numeric-looking
# strings (e.g. "268.91", "43.772.63") as text instead of auto-converting them
# to numbers, while leaving the cell style unchanged (General / default, no
# explicit number format) just like the source file.
function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue "D2" '43.772.63'
$ws.Range("E2").Value = '  -0.43%  '
Set-TextValue "D3" '2.300.14'
$ws.Range("E3").Value = '  +3.60%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue "D5" '268.91'
$ws.Range("E5").Value = '  +2.15%  '
Set-TextValue "D6" '94.47'
$ws.Range("E6").Value = '  +8.36%  '
Set-TextValue "D7" '0.627'
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("E8").Value = '  +0.03%  '
Set-TextValue "D9" '0.624'
$ws.Range("E9").Value = '  +2.66%  '
Set-TextValue "D10" '45.86'
$ws.Range("E10").Value = '  +0.41%  '
Set-TextValue "D11" '0.0936'
$ws.Range("E11").Value = '  +1.47%  '
Set-TextValue "D12" '8.06'
$ws.Range("E12").Value = '  +6.34%  '
$ws.Range("E13").Value = '  +0.32%  '
Set-TextValue "D14" '2.641.27'
$ws.Range("E14").Value = '  +3.49%  '
Set-TextValue "D15" '15.38'
$ws.Range("E15").Value = '  +4.95%  '
Set-TextValue "D16" '0.850'
$ws.Range("E16").Value = '  +7.88%  '
Set-TextValue "D17" '2.317.99'
$ws.Range("E17").Value = '  +4.92%  '
Set-TextValue "D18" '43.801.57'
$ws.Range("E18").Value = '  -0.27%  '
Set-TextValue "D19" '0.0000106'
$ws.Range("E19").Value = '  +1.46%  '
Set-TextValue "D20" '6.31'
$ws.Range("E20").Value = '  +5.38%  '
Set-TextValue "D21" '71.15'
$ws.Range("E21").Value = '  +1.57%  '
$ws.Range("E22").Value = '  -3.34%  '
Set-TextValue "D23" '237.41'
$ws.Range("E23").Value = '  +2.42%  '
Set-TextValue "D24" '9.78'
$ws.Range("E24").Value = '  +8.90%  '
$ws.Range("E25").Value = '  +0.01%  '
Set-TextValue "D26" '11.28'
$ws.Range("E26").Value = '  +5.19%  '
Set-TextValue "D27" '2.49'
$ws.Range("E27").Value = '  +9.74%  '
$ws.Range("E28").Value = '  -3.78%  '
Set-TextValue "D29" '39.12'
$ws.Range("E29").Value = '  -1.56%  '
Set-TextValue "D30" '2.25'
$ws.Range("E30").Value = '  +1.38%  '
Set-TextValue "D31" '22.32'
$ws.Range("E31").Value = '  +8.48%  '
Set-TextValue "D32" '173.64'
$ws.Range("E32").Value = '  -0.71%  '
Set-TextValue "D33" '0.0897'
$ws.Range("E33").Value = '  +1.42%  '
Set-TextValue "D34" '5.55'
$ws.Range("E34").Value = '  +2.12%  '
$ws.Range("E35").Value = '  +1.46%  '
$ws.Range("E36").Value = '  -0.27%  '
Set-TextValue "D37" '4.50'
$ws.Range("E37").Value = '  +0.56%  '
Set-TextValue "D38" '0.0351'
$ws.Range("E38").Value = '  -2.57%  '
Set-TextValue "D39" '3.42'
$ws.Range("E39").Value = '  +4.89%  '
$ws.Range("E40").Value = '  +16.54%  '
$ws.Range("E42").Value = '  -0.96%  '
Set-TextValue "D43" '1.33'
$ws.Range("E43").Value = '  +17.46%  '
$ws.Range("E44").Value = '  -1.62%  '
Set-TextValue "D45" '61.59'
$ws.Range("E45").Value = '  -5.67%  '
Set-TextValue "D46" '8.85'
$ws.Range("E46").Value = '  +5.81%  '
$ws.Range("E47").Value = '  +3.75%  '
Set-TextValue "D48" '100.50'
$ws.Range("E48").Value = '  -0.86%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue "D50" '2.522.47'
$ws.Range("E50").Value = '  +3.60%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue "D51" '0.430'
$ws.Range("E51").Value = '  -4.18%  '
